$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C to FALSE for the existing 442 rows
$ws.Range("C1:C442").Value = $false

# Append three new rows with new vacancy IDs
$newRows = @(
    @{ Id = "72801067"; Url = "https://api.hh.ru/vacancies/72801067?host=hh.ru " },
    @{ Id = "72798766"; Url = "https://api.hh.ru/vacancies/72798766?host=hh.ru " },
    @{ Id = "72784352"; Url = "https://api.hh.ru/vacancies/72784352?host=hh.ru " }
)

$startRow = 443
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $newRows[$i].Id
    $ws.Cells.Item($r, 2).Value = $newRows[$i].Url
    $ws.Cells.Item($r, 3).Value = $false
}
